# Add the new daily profit row for 09/20/2025, matching the existing
# layout where column A holds the date as literal text and column B
# holds the profit as a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A to Text so the "MM/DD/YYYY" string is stored literally
# instead of being auto-converted into a date serial number, then reset
# the cell style back to Normal so no stray formatting is left behind.
$ws.Range("A34").NumberFormat = "@"
$ws.Range("A34").Value = "09/20/2025"
$ws.Range("A34").Style = "Normal"

$ws.Range("B34").Value = 15636.31
